$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 509:510 (existing rows 509..559 shift down to 511..561)
$ws.Range("A509:A510").EntireRow.Insert()

# Row 509 - new weekly entry, date 2021-09-10 (serial 44449), "1a amarillo"
$ws.Range("A509").Value = 3
$ws.Range("B509").Value = "Femacal de La Calera"
$ws.Range("C509").Value = "Coquimbo"
$ws.Range("D509").Value = 44449
$ws.Range("E509").Value = 5
$ws.Range("F509").Value = "Fruta"
$ws.Range("G509").Value = 100102
$ws.Range("H509").Value = "Cítricos"
$ws.Range("I509").Value = 100102003
$ws.Range("J509").Value = "Limón"
$ws.Range("K509").Value = "Sin especificar"
$ws.Range("L509").Value = "1a amarillo"
$ws.Range("M509").Value = 332
$ws.Range("N509").Value = 3000
$ws.Range("O509").Value = 3500
$ws.Range("P509").Value = 3229
$ws.Range("Q509").Value = "`$/malla 16 kilos"
$ws.Range("R509").Value = "Provincia de Quillota"
$ws.Range("S509").Value = 202
$ws.Range("T509").Value = 16

# Row 510 - new weekly entry, date 2021-09-10 (serial 44449), "2a amarillo"
$ws.Range("A510").Value = 3
$ws.Range("B510").Value = "Femacal de La Calera"
$ws.Range("C510").Value = "Coquimbo"
$ws.Range("D510").Value = 44449
$ws.Range("E510").Value = 5
$ws.Range("F510").Value = "Fruta"
$ws.Range("G510").Value = 100102
$ws.Range("H510").Value = "Cítricos"
$ws.Range("I510").Value = 100102003
$ws.Range("J510").Value = "Limón"
$ws.Range("K510").Value = "Sin especificar"
$ws.Range("L510").Value = "2a amarillo"
$ws.Range("M510").Value = 220
$ws.Range("N510").Value = 2400
$ws.Range("O510").Value = 2500
$ws.Range("P510").Value = 2455
$ws.Range("Q510").Value = "`$/malla 16 kilos"
$ws.Range("R510").Value = "Provincia de Quillota"
$ws.Range("S510").Value = 153
$ws.Range("T510").Value = 16
